$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.922.13"
$ws.Range("E2").Value = "  -0.65%  "
$ws.Range("D3").Value = "2.045.72"
$ws.Range("E3").Value = "  -0.42%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "228.58"
$ws.Range("E5").Value = "  +0.14%  "
$ws.Range("E6").Value = "  -1.28%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "60.83"
$ws.Range("E7").Value = "  +0.19%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  -2.18%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0827"
$ws.Range("E10").Value = "  +0.00%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.105"
$ws.Range("E11").Value = "  +0.50%  "
$ws.Range("D12").Value = "2.348.62"
$ws.Range("E12").Value = "  -0.52%  "
$ws.Range("E13").Value = "  -0.96%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.10"
$ws.Range("E14").Value = "  +0.18%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.779"
$ws.Range("E15").Value = "  +2.33%  "
$ws.Range("E16").Value = "  -1.21%  "
$ws.Range("D17").Value = "2.054.14"
$ws.Range("E17").Value = "  -0.05%  "
$ws.Range("D18").Value = "37.880.48"
$ws.Range("E18").Value = "  -0.62%  "
$ws.Range("E19").Value = "  -0.36%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.94"
$ws.Range("E20").Value = "  -3.82%  "
$ws.Range("E21").Value = "  -1.00%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "224.19"
$ws.Range("E22").Value = "  -0.36%  "
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("E24").Value = "  -0.19%  "
$ws.Range("E25").Value = "  +2.97%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "168.15"
$ws.Range("E26").Value = "  +1.18%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.33"
$ws.Range("E27").Value = "  +1.20%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.130"
$ws.Range("E28").Value = "  -1.69%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.80"
$ws.Range("E29").Value = "  -0.21%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.27"
$ws.Range("E30").Value = "  -1.95%  "
$ws.Range("E31").Value = "  -0.49%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.21"
$ws.Range("E32").Value = "  +7.62%  "
$ws.Range("E33").Value = "  -2.01%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.52"
$ws.Range("E34").Value = "  -0.76%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0604"
$ws.Range("E35").Value = "  -0.31%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.57"
$ws.Range("E36").Value = "  +4.48%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.35"
$ws.Range("E37").Value = "  +2.97%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.45"
$ws.Range("E38").Value = "  +6.01%  "
$ws.Range("E39").Value = "  -0.04%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.14"
$ws.Range("E40").Value = "  +8.10%  "
$ws.Range("D41").Value = "1.536.98"
$ws.Range("E41").Value = "  +0.57%  "
$ws.Range("E42").Value = "  -1.17%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "96.61"
$ws.Range("E43").Value = "  -1.47%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.82"
$ws.Range("E44").Value = "  -0.54%  "
$ws.Range("E45").Value = "  -1.65%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.15"
$ws.Range("E46").Value = "  +4.14%  "
$ws.Range("E47").Value = "  -0.75%  "
$ws.Range("E48").Value = "  -0.59%  "
$ws.Range("E49").Value = "  -1.36%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.09"
$ws.Range("E50").Value = "  +0.36%  "
$ws.Range("D51").Value = "2.238.52"
$ws.Range("E51").Value = "  -0.49%  "
